$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "The interest Rates has a format error"
$ws.Range("G4").Value = "2022-09-06 12:53:09"
